# Scheduled market-data refresh: update Leve price/profit columns (H-N)
# across all craft-log sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 883.625
$ws.Range("J17").Value = 883.9836
$ws.Range("L17").Value = 2651.9508
$ws.Range("N17").Value = -2987.9508

$ws.Range("H33").Value = 226.1875
$ws.Range("I33").Value = 237.61539
$ws.Range("K33").Value = 237.61539
$ws.Range("M33").Value = -8.615389999999991

$ws.Range("H64").Value = 7069.3335
$ws.Range("I64").Value = 3676.6667
$ws.Range("K64").Value = 3676.6667
$ws.Range("M64").Value = -3428.6667

$ws.Range("H67").Value = 7069.3335
$ws.Range("I67").Value = 3676.6667
$ws.Range("K67").Value = 3676.6667
$ws.Range("M67").Value = -2818.6667

$ws.Range("H112").Value = 1109.8379
$ws.Range("J112").Value = 1125.4412
$ws.Range("L112").Value = 3376.3236
$ws.Range("N112").Value = -5592.3236

$ws.Range("H129").Value = 1998.1666
$ws.Range("I129").Value = 1797.8
$ws.Range("K129").Value = 5393.4
$ws.Range("M129").Value = -393.3999999999996

$ws.Range("H132").Value = 36732.33
$ws.Range("I132").Value = 42829.805
$ws.Range("K132").Value = 128489.415
$ws.Range("M132").Value = -125959.415

$ws.Range("H137").Value = 32265814
$ws.Range("I137").Value = 83334340
$ws.Range("J137").Value = 1624698.5
$ws.Range("K137").Value = 250003020
$ws.Range("L137").Value = 4874095.5
$ws.Range("M137").Value = -250000470
$ws.Range("N137").Value = -4879195.5

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H141").Value = 1323.2
$ws.Range("I141").Value = 1323.2
$ws.Range("K141").Value = 3969.6
$ws.Range("M141").Value = 1210.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7364.375
$ws.Range("I2").Value = 6486.3335
$ws.Range("K2").Value = 6486.3335
$ws.Range("M2").Value = -6373.3335

$ws.Range("H32").Value = 6667579
$ws.Range("I32").Value = 6850173
$ws.Range("J32").Value = 2900
$ws.Range("K32").Value = 6850173
$ws.Range("L32").Value = 2900
$ws.Range("M32").Value = -6849886
$ws.Range("N32").Value = -3474

$ws.Range("H45").Value = 2209.476
$ws.Range("I45").Value = 2073.6843
$ws.Range("J45").Value = 3499.5
$ws.Range("K45").Value = 2073.6843
$ws.Range("L45").Value = 3499.5
$ws.Range("M45").Value = -1696.6843
$ws.Range("N45").Value = -4253.5

$ws.Range("H74").Value = 3206852
$ws.Range("I74").Value = 4465264
$ws.Range("J74").Value = 3621.818
$ws.Range("K74").Value = 4465264
$ws.Range("L74").Value = 3621.818
$ws.Range("M74").Value = -4464390
$ws.Range("N74").Value = -5369.818

$ws.Range("H77").Value = 3206852
$ws.Range("I77").Value = 4465264
$ws.Range("J77").Value = 3621.818
$ws.Range("K77").Value = 22326320
$ws.Range("L77").Value = 18109.09
$ws.Range("M77").Value = -22321952
$ws.Range("N77").Value = -26845.09

$ws.Range("H88").Value = 1449.6666
$ws.Range("I88").Value = 1399.6666
$ws.Range("J88").Value = 1499.6666
$ws.Range("K88").Value = 1399.6666
$ws.Range("L88").Value = 1499.6666
$ws.Range("M88").Value = -993.6666
$ws.Range("N88").Value = -2311.6666

$ws.Range("H91").Value = 1449.6666
$ws.Range("I91").Value = 1399.6666
$ws.Range("J91").Value = 1499.6666
$ws.Range("K91").Value = 1399.6666
$ws.Range("L91").Value = 1499.6666
$ws.Range("M91").Value = 4.333399999999983
$ws.Range("N91").Value = -4307.6666

$ws.Range("H110").Value = 1667.7273
$ws.Range("I110").Value = 1640.2
$ws.Range("K110").Value = 1640.2
$ws.Range("M110").Value = 404.8

$ws.Range("H116").Value = 7364.375
$ws.Range("I116").Value = 6486.3335
$ws.Range("K116").Value = 6486.3335
$ws.Range("M116").Value = -4192.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7364.375
$ws.Range("I3").Value = 6486.3335
$ws.Range("K3").Value = 6486.3335
$ws.Range("M3").Value = -6372.3335

$ws.Range("H20").Value = 1628.4546
$ws.Range("I20").Value = 1247.3334
$ws.Range("K20").Value = 1247.3334
$ws.Range("M20").Value = -1000.3334

$ws.Range("H54").Value = 4362.375
$ws.Range("I54").Value = 3316.6667
$ws.Range("J54").Value = 7499.5
$ws.Range("K54").Value = 3316.6667
$ws.Range("L54").Value = 7499.5
$ws.Range("M54").Value = -2832.6667
$ws.Range("N54").Value = -8467.5

$ws.Range("H105").Value = 274.33334
$ws.Range("I105").Value = 274.33334
$ws.Range("K105").Value = 274.33334
$ws.Range("M105").Value = 1472.66666

$ws.Range("H134").Value = 686905.4399999999
$ws.Range("I134").Value = 757133.9
$ws.Range("J134").Value = 462174.4
$ws.Range("K134").Value = 2271401.7
$ws.Range("L134").Value = 1386523.2
$ws.Range("M134").Value = -2268866.7
$ws.Range("N134").Value = -1391593.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 365323.06
$ws.Range("I58").Value = 441946
$ws.Range("J58").Value = 7749.3335
$ws.Range("K58").Value = 441946
$ws.Range("L58").Value = 7749.3335
$ws.Range("M58").Value = -441743
$ws.Range("N58").Value = -8155.3335

$ws.Range("H136").Value = 365323.06
$ws.Range("I136").Value = 441946
$ws.Range("J136").Value = 7749.3335
$ws.Range("K136").Value = 1325838
$ws.Range("L136").Value = 23248.0005
$ws.Range("M136").Value = -1323288
$ws.Range("N136").Value = -28348.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 435.77777
$ws.Range("I8").Value = 435.77777
$ws.Range("K8").Value = 1307.33331
$ws.Range("M8").Value = -1168.33331

$ws.Range("H133").Value = 7691.1113
$ws.Range("I133").Value = 5635.857
$ws.Range("K133").Value = 16907.571
$ws.Range("M133").Value = -11847.571

$ws.Range("H137").Value = 4200.5
$ws.Range("J137").Value = 5273.5454
$ws.Range("L137").Value = 15820.6362
$ws.Range("N137").Value = -26020.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5978.5454
$ws.Range("I70").Value = 5519
$ws.Range("J70").Value = 6530
$ws.Range("K70").Value = 5519
$ws.Range("L70").Value = 6530
$ws.Range("M70").Value = -5249
$ws.Range("N70").Value = -7070

$ws.Range("H73").Value = 5978.5454
$ws.Range("I73").Value = 5519
$ws.Range("J73").Value = 6530
$ws.Range("K73").Value = 5519
$ws.Range("L73").Value = 6530
$ws.Range("M73").Value = -4583
$ws.Range("N73").Value = -8402

$ws.Range("H102").Value = 1647.625
$ws.Range("I102").Value = 1731.8334
$ws.Range("J102").Value = 1395
$ws.Range("K102").Value = 1731.8334
$ws.Range("L102").Value = 1395
$ws.Range("M102").Value = -109.8334
$ws.Range("N102").Value = -4639

$ws.Range("H113").Value = 3484.875
$ws.Range("I113").Value = 3247.5
$ws.Range("K113").Value = 3247.5
$ws.Range("M113").Value = -1077.5

$ws.Range("H132").Value = 13152802
$ws.Range("I132").Value = 17764402
$ws.Range("J132").Value = 9741.85
$ws.Range("K132").Value = 53293206
$ws.Range("L132").Value = 29225.55
$ws.Range("M132").Value = -53290676
$ws.Range("N132").Value = -34285.55

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 447.875
$ws.Range("I16").Value = 434.14285
$ws.Range("K16").Value = 434.14285
$ws.Range("M16").Value = -264.14285

$ws.Range("H40").Value = 3766.5881
$ws.Range("I40").Value = 4040.923
$ws.Range("K40").Value = 4040.923
$ws.Range("M40").Value = -3904.923

$ws.Range("H61").Value = 2166.3333
$ws.Range("J61").Value = 2499.5
$ws.Range("L61").Value = 2499.5
$ws.Range("N61").Value = -2903.5

$ws.Range("H113").Value = 2166.3333
$ws.Range("J113").Value = 2499.5
$ws.Range("L113").Value = 2499.5
$ws.Range("N113").Value = -6839.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7313.357
$ws.Range("I81").Value = 7041.4287
$ws.Range("K81").Value = 14082.8574
$ws.Range("M81").Value = -13021.8574

$ws.Range("H84").Value = 7313.357
$ws.Range("I84").Value = 7041.4287
$ws.Range("K84").Value = 70414.28700000001
$ws.Range("M84").Value = -65110.28700000001

$ws.Range("H132").Value = 6580075
$ws.Range("I132").Value = 8360659.5
$ws.Range("J132").Value = 5609
$ws.Range("K132").Value = 25081978.5
$ws.Range("L132").Value = 16827
$ws.Range("M132").Value = -25079448.5
$ws.Range("N132").Value = -21887

$ws.Range("H135").Value = 98735.836
$ws.Range("J135").Value = 98735.836
$ws.Range("L135").Value = 98735.836
$ws.Range("N135").Value = -108875.836
